# Auto-generated Excel COM-interop edit script
# Applies numeric corrections to the Leve profit calculation sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) produced by the scheduled
# market-data refresh. Each row updates currentAveragePrice / Leve cost
# columns (H-L) and the recomputed profit columns (M, N); some profit
# cells are cleared or newly populated depending on sign changes.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 75
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
# Row 78
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
# Row 129
$ws.Range("H129").Value = 1395.4634
$ws.Range("I129").Value = 606.1
$ws.Range("J129").Value = 1650.0968
$ws.Range("K129").Value = 1818.3
$ws.Range("L129").Value = 4950.2904
$ws.Range("M129").Value = 3181.7
$ws.Range("N129").Value = -14950.2904
# Row 137
$ws.Range("H137").Value = 2094.0688
$ws.Range("I137").Value = 3370.6667
$ws.Range("J137").Value = 1192.9412
$ws.Range("K137").Value = 10112.0001
$ws.Range("L137").Value = 3578.8236
$ws.Range("M137").Value = -7562.000100000001
$ws.Range("N137").Value = -8678.8236

$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 102.25
$ws.Range("I5").Value = 103
$ws.Range("K5").Value = 103
$ws.Range("M5").Value = 9
# Row 45
$ws.Range("H45").Value = 1861.7778
$ws.Range("I45").Value = 1608
$ws.Range("K45").Value = 1608
$ws.Range("M45").Value = -1231
# Row 97
$ws.Range("H97").Value = 23405.555
$ws.Range("I97").Value = 23405.555
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 23405.555
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -22909.555
$ws.Range("N97").ClearContents()
# Row 103
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
# Row 132
$ws.Range("H132").Value = 5025.744
$ws.Range("I132").Value = 4393
$ws.Range("J132").Value = 5217.485
$ws.Range("K132").Value = 13179
$ws.Range("L132").Value = 15652.455
$ws.Range("M132").Value = -10649
$ws.Range("N132").Value = -20712.455

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 102.25
$ws.Range("I4").Value = 103
$ws.Range("K4").Value = 103
$ws.Range("M4").Value = 12
# Row 97
$ws.Range("H97").Value = 1501
$ws.Range("I97").Value = 1501
$ws.Range("K97").Value = 1501
$ws.Range("M97").Value = -510
# Row 99
$ws.Range("H99").Value = 2445
$ws.Range("I99").Value = 1195
$ws.Range("J99").Value = 3695
$ws.Range("K99").Value = 1195
$ws.Range("L99").Value = 3695
$ws.Range("M99").Value = 303
$ws.Range("N99").Value = -6691
# Row 105
$ws.Range("H105").Value = 2651.6667
$ws.Range("I105").Value = 2333.6667
$ws.Range("J105").Value = 2969.6667
$ws.Range("K105").Value = 2333.6667
$ws.Range("L105").Value = 2969.6667
$ws.Range("M105").Value = -586.6667000000002
$ws.Range("N105").Value = -6463.6667
# Row 134
$ws.Range("H134").Value = 3589.7078
$ws.Range("I134").Value = 1587.359
$ws.Range("J134").Value = 6593.231
$ws.Range("K134").Value = 4762.076999999999
$ws.Range("L134").Value = 19779.693
$ws.Range("M134").Value = -2227.076999999999
$ws.Range("N134").Value = -24849.693

$ws = $wb.Worksheets.Item("CRP")
# Row 97
$ws.Range("H97").Value = 21001
$ws.Range("J97").Value = 21001
$ws.Range("L97").Value = 21001
$ws.Range("N97").Value = -22983
# Row 122
$ws.Range("H122").Value = 3543.5557
$ws.Range("I122").Value = 2315.3333
$ws.Range("J122").Value = 6000
$ws.Range("K122").Value = 6945.999899999999
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = -4495.999899999999
$ws.Range("N122").Value = -22900

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 189.6
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 189.6
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 568.8
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -914.8
# Row 80
$ws.Range("H80").Value = 84485.086
$ws.Range("I80").Value = 888
$ws.Range("J80").Value = 144197.28
$ws.Range("K80").Value = 2664
$ws.Range("L80").Value = 432591.84
$ws.Range("M80").Value = -1728
$ws.Range("N80").Value = -434463.84
# Row 83
$ws.Range("H83").Value = 84485.086
$ws.Range("I83").Value = 888
$ws.Range("J83").Value = 144197.28
$ws.Range("K83").Value = 7992
$ws.Range("L83").Value = 1297775.52
$ws.Range("M83").Value = -3312
$ws.Range("N83").Value = -1307135.52

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2648.5833
$ws.Range("I80").Value = 2580
$ws.Range("J80").Value = 2717.1667
$ws.Range("K80").Value = 2580
$ws.Range("L80").Value = 2717.1667
$ws.Range("M80").Value = -1582
$ws.Range("N80").Value = -4713.1667
# Row 83
$ws.Range("H83").Value = 2648.5833
$ws.Range("I83").Value = 2580
$ws.Range("J83").Value = 2717.1667
$ws.Range("K83").Value = 12900
$ws.Range("L83").Value = 13585.8335
$ws.Range("M83").Value = -7908
$ws.Range("N83").Value = -23569.8335
# Row 97
$ws.Range("H97").Value = 897
$ws.Range("I97").Value = 897
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 897
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -401
$ws.Range("N97").ClearContents()
# Row 132
$ws.Range("H132").Value = 2277.8975
$ws.Range("I132").Value = 1729.5
$ws.Range("K132").Value = 5188.5
$ws.Range("M132").Value = -2658.5

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 616.9375
$ws.Range("I22").Value = 416.2
$ws.Range("J22").Value = 708.1818
$ws.Range("K22").Value = 416.2
$ws.Range("L22").Value = 708.1818
$ws.Range("M22").Value = -121.2
$ws.Range("N22").Value = -1298.1818
# Row 27
$ws.Range("H27").Value = 616.9375
$ws.Range("I27").Value = 416.2
$ws.Range("J27").Value = 708.1818
$ws.Range("K27").Value = 416.2
$ws.Range("L27").Value = 708.1818
$ws.Range("M27").Value = -309.2
$ws.Range("N27").Value = -922.1818
# Row 46
$ws.Range("H46").Value = 625537.75
$ws.Range("I46").Value = 494
$ws.Range("J46").Value = 1000564
$ws.Range("K46").Value = 494
$ws.Range("L46").Value = 1000564
$ws.Range("M46").Value = -306
$ws.Range("N46").Value = -1000940
# Row 68
$ws.Range("H68").Value = 2250.2632
$ws.Range("I68").Value = 1396.25
$ws.Range("J68").Value = 3714.2856
$ws.Range("K68").Value = 1396.25
$ws.Range("L68").Value = 3714.2856
$ws.Range("M68").Value = -647.25
$ws.Range("N68").Value = -5212.2856
# Row 71
$ws.Range("H71").Value = 2250.2632
$ws.Range("I71").Value = 1396.25
$ws.Range("J71").Value = 3714.2856
$ws.Range("K71").Value = 6981.25
$ws.Range("L71").Value = 18571.428
$ws.Range("M71").Value = -3237.25
$ws.Range("N71").Value = -26059.428
# Row 76
$ws.Range("H76").Value = 5631.6665
$ws.Range("I76").Value = 997.5
$ws.Range("J76").Value = 14900
$ws.Range("K76").Value = 997.5
$ws.Range("L76").Value = 14900
$ws.Range("M76").Value = -659.5
$ws.Range("N76").Value = -15576
# Row 79
$ws.Range("H79").Value = 5631.6665
$ws.Range("I79").Value = 997.5
$ws.Range("J79").Value = 14900
$ws.Range("K79").Value = 997.5
$ws.Range("L79").Value = 14900
$ws.Range("M79").Value = 172.5
$ws.Range("N79").Value = -17240
# Row 82
$ws.Range("H82").Value = 1328.3043
$ws.Range("I82").Value = 1069.2727
$ws.Range("J82").Value = 1565.75
$ws.Range("K82").Value = 1069.2727
$ws.Range("L82").Value = 1565.75
$ws.Range("M82").Value = -708.2727
$ws.Range("N82").Value = -2287.75
# Row 85
$ws.Range("H85").Value = 1328.3043
$ws.Range("I85").Value = 1069.2727
$ws.Range("J85").Value = 1565.75
$ws.Range("K85").Value = 1069.2727
$ws.Range("L85").Value = 1565.75
$ws.Range("M85").Value = 178.7273
$ws.Range("N85").Value = -4061.75
# Row 100
$ws.Range("H100").Value = 3075.5
$ws.Range("I100").Value = 1350
$ws.Range("J100").Value = 3650.6667
$ws.Range("K100").Value = 1350
$ws.Range("L100").Value = 3650.6667
$ws.Range("M100").Value = -809
$ws.Range("N100").Value = -4732.6667

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 3332
$ws.Range("I62").Value = 3332
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3332
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2708
$ws.Range("N62").ClearContents()
# Row 65
$ws.Range("H65").Value = 3332
$ws.Range("I65").Value = 3332
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 16660
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -13540
$ws.Range("N65").ClearContents()
# Row 122
$ws.Range("H122").Value = 3222
$ws.Range("I122").Value = 1938.1111
$ws.Range("J122").Value = 8999.5
$ws.Range("K122").Value = 5814.3333
$ws.Range("L122").Value = 26998.5
$ws.Range("M122").Value = -3364.3333
$ws.Range("N122").Value = -31898.5
